$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.455.52'
$ws.Range('E2').Value = '  +0.41%  '

# Row 3
$ws.Range('D3').Value = '2.247.28'
$ws.Range('E3').Value = '  -0.12%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.16'
$ws.Range('E5').Value = '  -0.79%  '

# Row 6
$ws.Range('E6').Value = '  -0.16%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '75.66'
$ws.Range('E7').Value = '  -0.62%  '

# Row 8
$ws.Range('E8').Value = '  +0.10%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -1.83%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.05'
$ws.Range('E10').Value = '  +8.88%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  -0.46%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.27'
$ws.Range('E12').Value = '  +0.03%  '

# Row 13
$ws.Range('E13').Value = '  -1.42%  '

# Row 14
$ws.Range('D14').Value = '2.587.49'
$ws.Range('E14').Value = '  -0.07%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.57'
$ws.Range('E15').Value = '  -2.61%  '

# Row 16
$ws.Range('E16').Value = '  -1.07%  '

# Row 17
$ws.Range('D17').Value = '2.255.90'
$ws.Range('E17').Value = '  +0.77%  '

# Row 18
$ws.Range('D18').Value = '42.282.51'
$ws.Range('E18').Value = '  +0.09%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000103'
$ws.Range('E19').Value = '  +4.44%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').Value = '  -0.16%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.21'
$ws.Range('E21').Value = '  +0.74%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.24'
$ws.Range('E22').Value = '  +1.46%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '231.95'

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.26'
$ws.Range('E24').Value = '  +30.41%  '

# Row 25
$ws.Range('E25').Value = '  +0.03%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.50'
$ws.Range('E26').Value = '  +3.05%  '

# Row 27
$ws.Range('E27').Value = '  -3.35%  '

# Row 28
$ws.Range('E28').Value = '  -1.21%  '

# Row 29
$ws.Range('E29').Value = '  +1.47%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '167.95'
$ws.Range('E30').Value = '  -0.56%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.68'
$ws.Range('E31').Value = '  +0.41%  '

# Row 32
$ws.Range('E32').Value = '  -2.68%  '

# Row 33
$ws.Range('E33').Value = '  -0.73%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '30.56'
$ws.Range('E34').Value = '  -6.80%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.32'
$ws.Range('E35').Value = '  +11.04%  '

# Row 36
$ws.Range('E36').Value = '  -0.95%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.54'
$ws.Range('E37').Value = '  +0.58%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0316'
$ws.Range('E38').Value = '  +5.90%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '13.99'
$ws.Range('E39').Value = '  +5.74%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.18'
$ws.Range('E40').Value = '  -1.49%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.81'
$ws.Range('E41').Value = '  -2.89%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '64.04'
$ws.Range('E42').Value = '  +6.11%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.202'
$ws.Range('E43').Value = '  -0.72%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '108.16'
$ws.Range('E44').Value = '  -7.82%  '

# Row 45
$ws.Range('E45').Value = '  +0.69%  '

# Row 46
$ws.Range('E46').Value = '  +1.29%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.998'
$ws.Range('E47').Value = '  -0.03%  '

# Row 48
$ws.Range('E48').Value = '  -0.22%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.18'
$ws.Range('E49').Value = '  +0.53%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.34'
$ws.Range('E50').Value = '  +4.68%  '

# Row 51
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E51').Value = '  -4.07%  '
